$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Comparisons")

# Row 2 (Total)
$ws.Range("C2").Value = 25.393569946289062
$ws.Range("E2").Value = 2277.0

# Row 4 (Day 0)
$ws.Range("C4").Value = 6.831668376922607
$ws.Range("E4").Value = 450.0

# Row 6 (Day 2)
$ws.Range("C6").Value = 3.9416000843048096
$ws.Range("E6").Value = 377.0

# Row 8 (Day 4)
$ws.Range("C8").Value = 3.7400689125061035
$ws.Range("E8").Value = 362.0
